$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul",
    "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul", "Diane Paul",
    "Diane Paul", "Diane Paul", "Diane Paul", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell",
    "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell",
    "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Kenneth Howell", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller",
    "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller",
    "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Nancy Miller", "Thomas Clarke",
    "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke",
    "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke", "Thomas Clarke",
    "Thomas Clarke", "Thomas Clarke"
)

$values = @(
    45976, 45977, 45978, 45979, 45961, 45962, 45963, 45964, 45970, 45971,
    45972, 45973, 45988, 45989, 45990, 45991, 45993, 45994, 45995, 45996,
    45997, 45998, 45999, 46000, 45986, 45987, 45988, 45989, 46007, 46008,
    46009, 46010, 45959, 45960, 45961, 45962, 46012, 46013, 46014, 46015,
    46019, 46020, 46021, 46022, 45983, 45984, 45985, 45986, 45933, 45934,
    45935, 45936, 46018, 46019, 46020, 46021, 45969, 45970, 45971, 45972,
    45945, 45946, 45947, 45948, 45939, 45940, 45941, 45942, 45953, 45954,
    45955, 45956, 46017, 46018, 46019, 46020, 46005, 46006, 46007, 46008,
    46001, 46002, 46003, 46004, 45930, 45931, 45932, 45933, 45982, 45983,
    45984, 45985
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value2 = $values[$i]
}
